$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 612.3333
$ws.Range("I2").Value = 298.7143
$ws.Range("K2").Value = 298.7143
$ws.Range("M2").Value = -185.7143
$ws.Range("H19").Value = 1335.875
$ws.Range("I19").Value = 1300
$ws.Range("J19").Value = 1395.6666
$ws.Range("K19").Value = 1300
$ws.Range("L19").Value = 1395.6666
$ws.Range("M19").Value = -1125
$ws.Range("N19").Value = -1745.6666
$ws.Range("H33").Value = 236.33333
$ws.Range("I33").Value = 211
$ws.Range("J33").Value = 325
$ws.Range("K33").Value = 211
$ws.Range("L33").Value = 325
$ws.Range("M33").Value = 18
$ws.Range("N33").Value = -783
$ws.Range("H38").Value = 188.92857
$ws.Range("I38").Value = 188.92857
$ws.Range("K38").Value = 566.78571
$ws.Range("M38").Value = -194.78571
$ws.Range("H49").Value = 2756.5
$ws.Range("I49").Value = 2008.6666
$ws.Range("K49").Value = 6025.9998
$ws.Range("M49").Value = -5889.9998
$ws.Range("H53").Value = 392.66666
$ws.Range("I53").Value = 139
$ws.Range("J53").Value = 900
$ws.Range("K53").Value = 139
$ws.Range("L53").Value = 900
$ws.Range("M53").Value = 498
$ws.Range("N53").Value = -2174
$ws.Range("H58").Value = 666.75
$ws.Range("I58").Value = 666.75
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 2000.25
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -1850.25
$ws.Range("N58").ClearContents()
$ws.Range("H81").Value = 89999
$ws.Range("J81").Value = 89999
$ws.Range("L81").Value = 89999
$ws.Range("N81").Value = -91995
$ws.Range("H84").Value = 89999
$ws.Range("J84").Value = 89999
$ws.Range("L84").Value = 269997
$ws.Range("N84").Value = -279981
$ws.Range("H95").Value = 39521.89
$ws.Range("J95").Value = 39521.89
$ws.Range("L95").Value = 39521.89
$ws.Range("N95").Value = -45013.89
$ws.Range("H98").Value = 1380.2
$ws.Range("I98").Value = 1030.6522
$ws.Range("K98").Value = 1030.6522
$ws.Range("M98").Value = 467.3478
$ws.Range("H100").Value = 3770
$ws.Range("J100").Value = 3783.3333
$ws.Range("L100").Value = 3783.3333
$ws.Range("N100").Value = -4865.3333
$ws.Range("H115").Value = 2000
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = 2000
$ws.Range("K115").Value = 0
$ws.Range("L115").Value = 6000
$ws.Range("M115").ClearContents()
$ws.Range("N115").Value = -9134
$ws.Range("H122").Value = 1380.2
$ws.Range("I122").Value = 1030.6522
$ws.Range("K122").Value = 3091.9566
$ws.Range("M122").Value = -641.9566
$ws.Range("H125").Value = 8613.200000000001
$ws.Range("I125").Value = 5599.5
$ws.Range("J125").Value = 9076.846
$ws.Range("K125").Value = 50395.5
$ws.Range("L125").Value = 81691.614
$ws.Range("M125").Value = -47935.5
$ws.Range("N125").Value = -86611.614
$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -110120
$ws.Range("H135").Value = 55556010
$ws.Range("I135").Value = 62500468
$ws.Range("K135").Value = 562504212
$ws.Range("M135").Value = -562501677
$ws.Range("H137").Value = 3286.2964
$ws.Range("J137").Value = 3967.75
$ws.Range("L137").Value = 11903.25
$ws.Range("N137").Value = -17003.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2430.8235
$ws.Range("J2").Value = 3539.6
$ws.Range("L2").Value = 3539.6
$ws.Range("N2").Value = -3765.6
$ws.Range("H32").Value = 3246.102
$ws.Range("I32").Value = 2107.6597
$ws.Range("K32").Value = 2107.6597
$ws.Range("M32").Value = -1820.6597
$ws.Range("H61").Value = 90910820
$ws.Range("I61").Value = 125001310
$ws.Range("J61").Value = 2832.3333
$ws.Range("K61").Value = 125001310
$ws.Range("L61").Value = 2832.3333
$ws.Range("M61").Value = -125001098
$ws.Range("N61").Value = -3256.3333
$ws.Range("H62").Value = 65000
$ws.Range("J62").Value = 65000
$ws.Range("L62").Value = 65000
$ws.Range("N62").Value = -66248
$ws.Range("H65").Value = 65000
$ws.Range("J65").Value = 65000
$ws.Range("L65").Value = 195000
$ws.Range("N65").Value = -201240
$ws.Range("H74").Value = 29416804
$ws.Range("I74").Value = 35719428
$ws.Range("J74").Value = 4562.3335
$ws.Range("K74").Value = 35719428
$ws.Range("L74").Value = 4562.3335
$ws.Range("M74").Value = -35718554
$ws.Range("N74").Value = -6310.3335
$ws.Range("H77").Value = 29416804
$ws.Range("I77").Value = 35719428
$ws.Range("J77").Value = 4562.3335
$ws.Range("K77").Value = 178597140
$ws.Range("L77").Value = 22811.6675
$ws.Range("M77").Value = -178592772
$ws.Range("N77").Value = -31547.6675
$ws.Range("H110").Value = 67856.53
$ws.Range("I110").Value = 78138.38
$ws.Range("K110").Value = 78138.38
$ws.Range("M110").Value = -76093.38
$ws.Range("H116").Value = 2430.8235
$ws.Range("J116").Value = 3539.6
$ws.Range("L116").Value = 3539.6
$ws.Range("N116").Value = -8127.6
$ws.Range("H132").Value = 2781565
$ws.Range("I132").Value = 3033734.5
$ws.Range("J132").Value = 7700
$ws.Range("K132").Value = 9101203.5
$ws.Range("L132").Value = 23100
$ws.Range("M132").Value = -9098673.5
$ws.Range("N132").Value = -28160
$ws.Range("H136").Value = 90910820
$ws.Range("I136").Value = 125001310
$ws.Range("J136").Value = 2832.3333
$ws.Range("K136").Value = 375003930
$ws.Range("L136").Value = 8496.999899999999
$ws.Range("M136").Value = -375001380
$ws.Range("N136").Value = -13596.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2430.8235
$ws.Range("J3").Value = 3539.6
$ws.Range("L3").Value = 3539.6
$ws.Range("N3").Value = -3767.6
$ws.Range("H22").Value = 2357.4167
$ws.Range("I22").Value = 3082.3333
$ws.Range("J22").Value = 1632.5
$ws.Range("K22").Value = 3082.3333
$ws.Range("L22").Value = 1632.5
$ws.Range("M22").Value = -2909.3333
$ws.Range("N22").Value = -1978.5
$ws.Range("H86").Value = 4051.5
$ws.Range("I86").Value = 4100
$ws.Range("K86").Value = 4100
$ws.Range("M86").Value = -2977
$ws.Range("H89").Value = 4051.5
$ws.Range("I89").Value = 4100
$ws.Range("K89").Value = 20500
$ws.Range("M89").Value = -14884
$ws.Range("H94").Value = 1944.6842
$ws.Range("I94").Value = 2277.7273
$ws.Range("J94").Value = 1486.75
$ws.Range("K94").Value = 2277.7273
$ws.Range("L94").Value = 1486.75
$ws.Range("M94").Value = -1826.7273
$ws.Range("N94").Value = -2388.75
$ws.Range("H99").Value = 1565.7778
$ws.Range("I99").Value = 1635.125
$ws.Range("K99").Value = 1635.125
$ws.Range("M99").Value = -137.125
$ws.Range("H107").Value = 67353.56
$ws.Range("I107").Value = 4780.6665
$ws.Range("J107").Value = 255072.25
$ws.Range("K107").Value = 4780.6665
$ws.Range("L107").Value = 255072.25
$ws.Range("M107").Value = -2860.6665
$ws.Range("N107").Value = -258912.25
$ws.Range("H134").Value = 20003698
$ws.Range("I134").Value = 20836816
$ws.Range("J134").Value = 8897
$ws.Range("K134").Value = 62510448
$ws.Range("L134").Value = 26691
$ws.Range("M134").Value = -62507913
$ws.Range("N134").Value = -31761

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 456.6
$ws.Range("I22").Value = 445.5
$ws.Range("J22").Value = 464
$ws.Range("K22").Value = 445.5
$ws.Range("L22").Value = 464
$ws.Range("M22").Value = -95.5
$ws.Range("N22").Value = -1164
$ws.Range("H31").Value = 2841.4666
$ws.Range("I31").Value = 3133.4546
$ws.Range("J31").Value = 2038.5
$ws.Range("K31").Value = 3133.4546
$ws.Range("L31").Value = 2038.5
$ws.Range("M31").Value = -2838.4546
$ws.Range("N31").Value = -2628.5
$ws.Range("H34").Value = 2841.4666
$ws.Range("I34").Value = 3133.4546
$ws.Range("J34").Value = 2038.5
$ws.Range("K34").Value = 3133.4546
$ws.Range("L34").Value = 2038.5
$ws.Range("M34").Value = -2931.4546
$ws.Range("N34").Value = -2442.5
$ws.Range("H58").Value = 23815736
$ws.Range("J58").Value = 1773.8889
$ws.Range("L58").Value = 1773.8889
$ws.Range("N58").Value = -2179.8889
$ws.Range("H62").Value = 99
$ws.Range("J62").Value = 99
$ws.Range("L62").Value = 99
$ws.Range("N62").Value = -1347
$ws.Range("H65").Value = 99
$ws.Range("J65").Value = 99
$ws.Range("L65").Value = 495
$ws.Range("N65").Value = -6735
$ws.Range("H86").Value = 11698.8
$ws.Range("I86").Value = 8860
$ws.Range("K86").Value = 8860
$ws.Range("M86").Value = -7737
$ws.Range("H89").Value = 11698.8
$ws.Range("I89").Value = 8860
$ws.Range("K89").Value = 44300
$ws.Range("M89").Value = -38684
$ws.Range("H94").Value = 1670.0714
$ws.Range("I94").Value = 1925.4445
$ws.Range("J94").Value = 1210.4
$ws.Range("K94").Value = 1925.4445
$ws.Range("L94").Value = 1210.4
$ws.Range("M94").Value = -1474.4445
$ws.Range("N94").Value = -2112.4
$ws.Range("H109").Value = 21999
$ws.Range("J109").Value = 21999
$ws.Range("L109").Value = 21999
$ws.Range("N109").Value = -24079
$ws.Range("H132").Value = 71431096
$ws.Range("I132").Value = 76925530
$ws.Range("K132").Value = 230776590
$ws.Range("M132").Value = -230774060
$ws.Range("H134").Value = 9302238
$ws.Range("I134").Value = 10045879
$ws.Range("K134").Value = 30137637
$ws.Range("M134").Value = -30135102
$ws.Range("H136").Value = 23815736
$ws.Range("J136").Value = 1773.8889
$ws.Range("L136").Value = 5321.6667
$ws.Range("N136").Value = -10421.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 6840.091
$ws.Range("I3").Value = 5026.8887
$ws.Range("K3").Value = 15080.6661
$ws.Range("M3").Value = -14968.6661
$ws.Range("H32").Value = 250747.75
$ws.Range("I32").Value = 250747.75
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 752243.25
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -751960.25
$ws.Range("N32").ClearContents()
$ws.Range("H49").Value = 459.8
$ws.Range("I49").Value = 324.75
$ws.Range("K49").Value = 974.25
$ws.Range("M49").Value = -818.25
$ws.Range("H70").Value = 12405.272
$ws.Range("I70").Value = 3576.3333
$ws.Range("K70").Value = 10728.9999
$ws.Range("M70").Value = -10413.9999
$ws.Range("H73").Value = 12405.272
$ws.Range("I73").Value = 3576.3333
$ws.Range("K73").Value = 10728.9999
$ws.Range("M73").Value = -9636.999899999999
$ws.Range("H86").Value = 752.5625
$ws.Range("J86").Value = 839.8889
$ws.Range("L86").Value = 2519.6667
$ws.Range("N86").Value = -4891.6667
$ws.Range("H89").Value = 752.5625
$ws.Range("J89").Value = 839.8889
$ws.Range("L89").Value = 7559.0001
$ws.Range("N89").Value = -19415.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2452.3333
$ws.Range("I80").Value = 2832
$ws.Range("K80").Value = 2832
$ws.Range("M80").Value = -1834
$ws.Range("H83").Value = 2452.3333
$ws.Range("I83").Value = 2832
$ws.Range("K83").Value = 14160
$ws.Range("M83").Value = -9168
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492
$ws.Range("H102").Value = 3737.0908
$ws.Range("I102").Value = 3910.9
$ws.Range("J102").Value = 1999
$ws.Range("K102").Value = 3910.9
$ws.Range("L102").Value = 1999
$ws.Range("M102").Value = -2288.9
$ws.Range("N102").Value = -5243
$ws.Range("H107").Value = 4283.6665
$ws.Range("I107").Value = 4209.778
$ws.Range("J107").Value = 4357.5557
$ws.Range("K107").Value = 4209.778
$ws.Range("L107").Value = 4357.5557
$ws.Range("M107").Value = -2289.778
$ws.Range("N107").Value = -8197.555700000001
$ws.Range("H113").Value = 76213.28999999999
$ws.Range("I113").Value = 88399.664
$ws.Range("K113").Value = 88399.664
$ws.Range("M113").Value = -86229.664
$ws.Range("H122").Value = 78469.94
$ws.Range("I122").Value = 94501.46000000001
$ws.Range("K122").Value = 283504.38
$ws.Range("M122").Value = -281054.38
$ws.Range("H126").Value = 3640
$ws.Range("J126").Value = 8500
$ws.Range("L126").Value = 25500
$ws.Range("N126").Value = -30440
$ws.Range("H132").Value = 6584381
$ws.Range("I132").Value = 8338296
$ws.Range("K132").Value = 25014888
$ws.Range("M132").Value = -25012358

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3200
$ws.Range("I22").Value = 4230
$ws.Range("K22").Value = 4230
$ws.Range("M22").Value = -3935
$ws.Range("H27").Value = 3200
$ws.Range("I27").Value = 4230
$ws.Range("K27").Value = 4230
$ws.Range("M27").Value = -4123
$ws.Range("H40").Value = 3223.1333
$ws.Range("I40").Value = 3203.3572
$ws.Range("K40").Value = 3203.3572
$ws.Range("M40").Value = -3067.3572
$ws.Range("I68").Value = 13158894
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 13158894
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -13158145
$ws.Range("N68").ClearContents()
$ws.Range("I71").Value = 13158894
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 65794470
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -65790726
$ws.Range("N71").ClearContents()
$ws.Range("H93").Value = 3180.7
$ws.Range("I93").Value = 1149.6
$ws.Range("K93").Value = 1149.6
$ws.Range("M93").Value = 98.40000000000009
$ws.Range("H132").Value = 10008913
$ws.Range("I132").Value = 10008913
$ws.Range("K132").Value = 30026739
$ws.Range("M132").Value = -30024209
$ws.Range("H136").Value = 2699
$ws.Range("J136").Value = 2898
$ws.Range("L136").Value = 8694
$ws.Range("N136").Value = -13794

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 1500
$ws.Range("I33").Value = 2000
$ws.Range("J33").Value = 1000
$ws.Range("K33").Value = 2000
$ws.Range("L33").Value = 1000
$ws.Range("M33").Value = -1750
$ws.Range("N33").Value = -1500
$ws.Range("H36").Value = 1500
$ws.Range("I36").Value = 2000
$ws.Range("J36").Value = 1000
$ws.Range("K36").Value = 2000
$ws.Range("L36").Value = 1000
$ws.Range("M36").Value = -1750
$ws.Range("N36").Value = -1500
$ws.Range("H41").Value = 34241.5
$ws.Range("J41").Value = 34241.5
$ws.Range("L41").Value = 34241.5
$ws.Range("N41").Value = -35021.5
$ws.Range("H42").Value = 100000
$ws.Range("J42").Value = 100000
$ws.Range("L42").Value = 100000
$ws.Range("N42").Value = -100756
$ws.Range("H43").Value = 31966.666
$ws.Range("I43").Value = 31966.666
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 31966.666
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -31817.666
$ws.Range("N43").ClearContents()
$ws.Range("H100").Value = 895.2
$ws.Range("I100").Value = 835.04346
$ws.Range("K100").Value = 1670.08692
$ws.Range("M100").Value = -1129.08692
$ws.Range("H101").Value = 95149.25
$ws.Range("I101").Value = 180000
$ws.Range("J101").Value = 66865.664
$ws.Range("K101").Value = 180000
$ws.Range("L101").Value = 66865.664
$ws.Range("M101").Value = -176755
$ws.Range("N101").Value = -73355.664
$ws.Range("H122").Value = 3562.6667
$ws.Range("I122").Value = 3323.2856
$ws.Range("J122").Value = 3897.8
$ws.Range("K122").Value = 9969.856800000001
$ws.Range("L122").Value = 11693.4
$ws.Range("M122").Value = -7519.856800000001
$ws.Range("N122").Value = -16593.4
